# [REFIX] resources/users- password back from hash-> 12345
#
# Replace the bcrypt password-hash values in column D (rows 2-13) of the
# "Users" sheet with the plain numeric placeholder password 12345, and
# leave the selection on D2:D13 (active cell D2), matching the state Excel
# ends up in after retyping the values over that range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = 12345
    # Touch the font so Excel records an explicit style for the
    # retyped numeric cells (mirrors the cellXfs entry added upstream).
    $cell.Font.ThemeColor = 1
}

$ws.Range("D2:D13").Select()
